# Add new Doc Type "Note" to the "Types" lookup sheet.
#
# The workbook has 4 sheets: Documents, Status, Priority, Types.
# "Types" (xl/worksheets/sheet4.xml) holds the list of values that feed the
# "Type" column's data-validation dropdown on the "Documents" sheet. We add
# one new row ("Note") right after the existing last entry ("Fact Sheet").

$wb = $excel.ActiveWorkbook

$wsTypes = $wb.Worksheets.Item("Types")
$wsDocs  = $wb.Worksheets.Item("Documents")

# --- Content change: append "Note" as a new available Type -----------------
# Row 13 / A13 currently holds the last entry ("Fact Sheet"); the new entry
# goes into row 14 / A14. Writing the value creates the shared-string entry
# and extends the sheet's used range (dimension) automatically.
$wsTypes.Range("A14").Value = "Note"

# --- View-state bookkeeping matching the authoring session -----------------
# While adding the row the author was on the Types sheet with C14 selected
# (and that sheet's zoom recorded at 100%); they then returned to the
# Documents sheet (which stays the active tab) with C14 selected there too.
$wsTypes.Activate()
$wsTypes.Range("C14").Select()
$excel.ActiveWindow.Zoom = 100

$wsDocs.Activate()
$wsDocs.Range("C14").Select()
